# Update countries & provincias Spain
# Applies the daily COVID data refresh: new totals for India (row 5) and
# Pakistan (row 21), Tailandia/Gambia swapping their (A-column) label order
# while Tailandia's figures get refreshed, the Islas Malvinas/Montserrat
# rows swapping their D/H figures (and thus their labels), and the
# "datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 07:08"

# --- Row 5: India ---
$ws.Range("B5").Value = 5400619
$ws.Range("C5").Value = 2389
$ws.Range("D5").Value = 4303043
$ws.Range("E5").Value = 1010802

# --- Row 21: Pakistan ---
$ws.Range("B21").Value = 305671
$ws.Range("C21").Value = 640
$ws.Range("D21").Value = 292303
$ws.Range("E21").Value = 6952
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 6416

# --- Rows 134/135: Tailandia & Gambia swap order, Tailandia data refreshed ---
$ws.Range("A134").Value = "Tailandia"
$ws.Range("B134").Value = 3506
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 3340
$ws.Range("E134").Value = 107
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 59

$ws.Range("A135").Value = "Gambia"
$ws.Range("B135").Value = 3504
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 1992
$ws.Range("E135").Value = 1404
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 108

# --- Rows 204/205: Santa Lucia & Timor Oriental swap order (values equal) ---
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# --- Rows 214/215: Montserrat & Islas Malvinas swap order + D/H swap ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
